$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text cells (coin names, links) ---
$textValues = @{
    'B14' = 'BitForexToken'
    'C14' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'B15' = 'TigerCash'
    'C15' = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
    'B16' = 'LEO'
    'C16' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    'B17' = 'BTSEToken'
    'C17' = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
    'B18' = 'BitpandaEcosystemToken'
    'C18' = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
    'B19' = 'MCDex'
    'C19' = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
    'B20' = 'ProBitToken'
    'C20' = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
    'B21' = 'ZBToken'
    'C21' = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
}
foreach ($addr in $textValues.Keys) {
    $ws.Range($addr).Value = $textValues[$addr]
}

# --- Numeric-looking text cells (price / volume %) ---
# These must stay TEXT (as in the source data), so force the cell to
# Text format before assigning, then clear the format again so no
# extra style is left behind on the cell.
$numericTextValues = @{
    'D2' = '306.15'
    'E2' = '-4.37%'
    'D3' = '39.91'
    'E3' = '-6.41%'
    'D4' = '5.037'
    'E4' = '-3.84%'
    'D5' = '0.07668'
    'E5' = '-5.87%'
    'D6' = '4.230'
    'E6' = '-2.34%'
    'D7' = '1.603'
    'E7' = '-8.44%'
    'D8' = '0.8853'
    'E8' = '-6.90%'
    'D9' = '0.1003'
    'E9' = '-9.74%'
    'D10' = '0.1727'
    'E10' = '-6.65%'
    'D11' = '0.08919'
    'E11' = '-4.20%'
    'D12' = '0.04390'
    'E12' = '-5.94%'
    'D13' = '0.1054'
    'E13' = '-0.38%'
    'D14' = '0.001273'
    'E14' = '-1.01%'
    'D15' = '0.005800'
    'E15' = '-4.38%'
    'D16' = '3.351'
    'E16' = '-0.41%'
    'D17' = '2.530'
    'E17' = '-0.30%'
    'D18' = '0.3361'
    'E18' = '-0.08%'
    'D19' = '7.079'
    'E19' = '-4.67%'
    'D20' = '0.1341'
    'E20' = '-3.78%'
    'D21' = '0.3204'
    'E21' = '19.29%'
    'D22' = '0.04211'
    'E22' = '0.47%'
    'D23' = '0.001194'
    'E23' = '-4.64%'
    'D24' = '0.004062'
    'E24' = '-5.45%'
    'D25' = '0.0001224'
    'E25' = '-6.10%'
    'E26' = '-0.16%'
    'D38' = '0.02337'
    'E38' = '-10.04%'
    'D39' = '0.05143'
    'E39' = '-6.28%'
    'D40' = '0.007975'
    'E40' = '2.02%'
    'D41' = '0.1319'
    'E41' = '-5.04%'
    'D42' = '0.006683'
    'E42' = '1.31%'
    'D43' = '0.001991'
    'E43' = '-6.34%'
    'D44' = '0.008475'
    'E44' = '-0.04%'
    'D45' = '0.3034'
    'E45' = '-12.42%'
    'D46' = '0.00006581'
    'E46' = '-5.63%'
    'D47' = '0.00000000752'
    'E47' = '0.04%'
    'D48' = '0.003400'
    'E48' = '-2.56%'
    'E49' = '41.36%'
    'D50' = '0.00002106'
    'E50' = '0.04%'
    'D51' = '0.0002006'
    'E51' = '0.04%'
}
foreach ($addr in $numericTextValues.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $numericTextValues[$addr]
    $cell.ClearFormats()
}
